$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-07 23:48:48"
$ws.Range("E3").Value = "2026-02-07 23:48:51"
$ws.Range("I3").Value = "0.3 mm"
$ws.Range("E4").Value = "2026-02-07 23:48:53"
$ws.Range("H4").Value = "'55%"
$ws.Range("E5").Value = "2026-02-07 23:48:56"
$ws.Range("E6").Value = "2026-02-07 23:48:58"
$ws.Range("E7").Value = "2026-02-07 23:49:00"
$ws.Range("H7").Value = "'50%"
$ws.Range("E8").Value = "2026-02-07 23:49:03"
$ws.Range("H8").Value = "'64%"
$ws.Range("E9").Value = "2026-02-07 23:49:06"
$ws.Range("O9").Value = "10.3 °C"
$ws.Range("E10").Value = "2026-02-07 23:49:08"
$ws.Range("E11").Value = "2026-02-07 23:49:11"
$ws.Range("O11").Value = "3.3 °C"
$ws.Range("E12").Value = "2026-02-07 23:49:13"
$ws.Range("H12").Value = "'85%"
$ws.Range("E13").Value = "2026-02-07 23:49:16"
$ws.Range("O13").Value = "2.9 °C"
$ws.Range("E14").Value = "2026-02-07 23:49:19"
$ws.Range("H14").Value = "'62%"
$ws.Range("I14").Value = "0.1 mm"
$ws.Range("O14").Value = "11.6 °C"
$ws.Range("E15").Value = "2026-02-07 23:49:22"
$ws.Range("H15").Value = "'72%"
$ws.Range("O15").Value = "10.0 °C"
$ws.Range("E16").Value = "2026-02-07 23:49:24"
$ws.Range("H16").Value = "'62%"
$ws.Range("I16").Value = "1.4 mm"
$ws.Range("E17").Value = "2026-02-07 23:49:27"
$ws.Range("E18").Value = "2026-02-07 23:49:30"
$ws.Range("E19").Value = "2026-02-07 23:49:32"
$ws.Range("E20").Value = "2026-02-07 23:49:35"
$ws.Range("I20").Value = "5.9 mm"
$ws.Range("E21").Value = "2026-02-07 23:49:37"
$ws.Range("I21").Value = "0.7 mm"
$ws.Range("J21").Value = "1006.2 hPa"
$ws.Range("O21").Value = "4.8 °C"
$ws.Range("E22").Value = "2026-02-07 23:49:40"
$ws.Range("I22").Value = "3.5 mm"
$ws.Range("E23").Value = "2026-02-07 23:49:43"
$ws.Range("E24").Value = "2026-02-07 23:49:45"
$ws.Range("J24").Value = "1007.0 hPa"
$ws.Range("E25").Value = "2026-02-07 23:49:48"
$ws.Range("E26").Value = "2026-02-07 23:49:51"
$ws.Range("L26").Value = "49.7 km/h - 229º 23:05 TU"
$ws.Range("E27").Value = "2026-02-07 23:49:54"
$ws.Range("I27").Value = "3.2 mm"
$ws.Range("E28").Value = "2026-02-07 23:49:57"
$ws.Range("H28").Value = "'64%"
$ws.Range("E29").Value = "2026-02-07 23:49:59"
$ws.Range("H29").Value = "'72%"
$ws.Range("O29").Value = "10.3 °C"
$ws.Range("E30").Value = "2026-02-07 23:50:02"
$ws.Range("E31").Value = "2026-02-07 23:50:05"
$ws.Range("H31").Value = "'68%"
$ws.Range("J31").Value = "1003.2 hPa"
$ws.Range("L31").Value = "77.8 km/h - 228º 23:29 TU"
$ws.Range("E32").Value = "2026-02-07 23:50:07"
$ws.Range("I32").Value = "2.3 mm"
$ws.Range("E33").Value = "2026-02-07 23:50:10"
$ws.Range("E34").Value = "2026-02-07 23:50:13"
$ws.Range("E35").Value = "2026-02-07 23:50:16"
$ws.Range("E36").Value = "2026-02-07 23:50:19"
$ws.Range("E37").Value = "2026-02-07 23:50:21"
$ws.Range("O37").Value = "4.6 °C"
$ws.Range("E38").Value = "2026-02-07 23:50:24"
$ws.Range("E39").Value = "2026-02-07 23:50:26"
$ws.Range("H39").Value = "'67%"
$ws.Range("O39").Value = "-4.6 °C"
$ws.Range("E40").Value = "2026-02-07 23:50:29"
$ws.Range("H40").Value = "'86%"
$ws.Range("I40").Value = "0.6 mm"
$ws.Range("E41").Value = "2026-02-07 23:50:32"
$ws.Range("I41").Value = "0.7 mm"
$ws.Range("O41").Value = "12.6 °C"
$ws.Range("E42").Value = "2026-02-07 23:50:35"
$ws.Range("E43").Value = "2026-02-07 23:50:37"
$ws.Range("E44").Value = "2026-02-07 23:50:40"
$ws.Range("E45").Value = "2026-02-07 23:50:43"
$ws.Range("E46").Value = "2026-02-07 23:50:45"
$ws.Range("J46").Value = "1007.2 hPa"
